# Fruta / hortaliza, semanal
# Insert a new weekly record for "Espárragos" at Feria Lagunitas de Puerto Montt.
# The new observation is inserted as row 40, pushing the existing rows 40-59
# down to 41-60 (dimension grows from A1:R59 to A1:R60).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 40..59 down to 41..60, opening up a blank row 40.
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new weekly data point.
$ws.Range("A40").Value = 4
$ws.Range("B40").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C40").Value = "Los Lagos"
$ws.Range("D40").Value = 44876
$ws.Range("E40").Value = 10
$ws.Range("F40").Value = 300000000
$ws.Range("G40").Value = "Espárragos"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 500
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1700
$ws.Range("M40").Value = 1600
$ws.Range("N40").Value = "$/kilo"
$ws.Range("O40").Value = "Provincia de Linares"
$ws.Range("P40").Value = 1600
$ws.Range("Q40").Value = 1
$ws.Range("R40").Value = "Hortaliza"
